# Finalized definition of virtual environment
# Update calibration numbers across AR / SETAR / GARCH / TARCH / AR_TARCH sheets

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# AR sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("AR")
$ws.Range("B2").Value = 0.0000008876168505136168
$ws.Range("B3").Value = 0.6547751660728359
$ws.Range("B4").Value = 0.000004066589104907773
$ws.Range("B5").Value = "[1.0, 0.43727518493587086, 0.28613477544328186, 0.27199442917878686, 0.3099611735917779, 0.467973256164935, 0.3396808570056091, 0.2448466431936482, 0.21821445137316098, 0.21377132039456312, 0.21902846430909068, 0.23691302729989808, 0.23541187232751404, 0.21071593789224602, 0.22419646267745036, 0.20948392104713245, 0.2303178288568498, 0.2303690946519094, 0.16204017765824769, 0.15069940842756444]"

# ---------------------------------------------------------------------------
# SETAR sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("SETAR")
$ws.Range("B3").Value = 1.051556420233463
$ws.Range("B4").Value = -0.0005949638109892461
$ws.Range("B5").Value = 0.6327209843211998
$ws.Range("B6").Value = 0.000002700084099889127
$ws.Range("B7").Value = 0.0006872867163584654
$ws.Range("B8").Value = 0.5987507056276733
$ws.Range("B9").Value = 0.000002676715403092161
$ws.Range("B10").Value = "[1.0, 0.38916709900314567, 0.2671351704386194, 0.2584173185811342, 0.20862369260020014, 0.3991024229982721, 0.28846610100682035, 0.20625907679339422, 0.21251455954133502, 0.21285034733319905, 0.18581937307552124, 0.22445793664741498, 0.1971510425511354, 0.17385409045206168, 0.225549725874473, 0.1984651964423956, 0.2124441653047936, 0.21902780591950471, 0.18081301793715276, 0.1549700575228758]"

# ---------------------------------------------------------------------------
# GARCH sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GARCH")
$ws.Range("B2").Value = -0.0000002751423149905367
$ws.Range("B3").Value = 0.00000009830199089580138
$ws.Range("B4").Value = 0.1
$ws.Range("B5").Value = 0.88
$ws.Range("B6").Value = "[1.0, 0.06395191982298506, -0.05351709397918391, -0.03223601164085689, 0.009227611377080277, 0.35157781501642676, 0.02263000124447862, -0.02276025665785845, -0.009789200777334141, -0.0002163434050001766, -0.023516358993410206, -0.049379264707421926, -7.515246758931471e-06, 0.004451247618535848, 0.012323789366523714, 0.003746767303173175, -0.02301872688112336, -0.0036016299964977054, -0.011621900899442782, -0.024313442163423897]"

# ---------------------------------------------------------------------------
# TARCH sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("TARCH")
$ws.Range("B2").Value = -0.0000002751423149905367
$ws.Range("B3").Value = 0.00000009830199089580138
$ws.Range("B4").Value = 0.1
$ws.Range("B5").Value = 0.875
$ws.Range("B6").Value = "[0.9999999999999999, 0.059842353250227055, -0.055751295046479646, -0.03380299768492689, 0.0066757047636034385, 0.3473572964904588, 0.020158682087921813, -0.02435268652801468, -0.01088036700055235, -0.001025194944400129, -0.0244682434237864, -0.05065735050401834, -0.0005925486518791323, 0.0039028803109815266, 0.012033765047920355, 0.0031947704540315003, -0.02382598553135098, -0.003421746699646458, -0.011418752972101669, -0.02434566323149332]"
$ws.Range("B7").Value = 0.01

# ---------------------------------------------------------------------------
# AR_TARCH sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("AR_TARCH")
$ws.Range("B2").Value = 0.0000008876168505136168
$ws.Range("B3").Value = 0.0000004066589104907772
$ws.Range("B5").Value = 0.6749999999999999
$ws.Range("B6").Value = "[1.0, 0.07920536197451054, 0.005793542288766162, 0.01207935709793957, 0.03393166865203931, 0.27318277173271677, 0.0447918692878711, -0.003285963069549128, 0.015354859408816578, 0.02847226749698569, 0.01855307658450718, 0.0237243130875621, 0.02757375004618389, 0.04510355106307968, 0.03369678535166264, 0.011182759189600286, 0.01637705077038864, 0.051800863366321614, 0.04368335472013011, -0.013865528431752063]"
$ws.Range("B7").Value = 0.05
$ws.Range("B9").Value = 0.6547751660728359
